$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top; existing data (rows 1-20) shifts down to rows 2-21
$ws.Rows.Item(1).Insert()

# New header row
$ws.Range("A1").Value = "Rank"
$ws.Range("B1").Value = "City Name"
$ws.Range("C1").Value = "Overnight International Visitors (Millions)"
$ws.Range("D1").Value = "Year"

# Fill the new Year column for each of the 20 data rows with 2012
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 4).Value = 2012
}

# Update selection to match the target state
$ws.Range("D2:D21").Select()
